$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 27092
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 672
$ws.Range("F6").Value = 585
$ws.Range("F8").Value = 387
$ws.Range("F9").Value = 505
$ws.Range("F11").Value = 56
$ws.Range("F12").Value = 321
$ws.Range("F13").Value = 108
$ws.Range("F14").Value = 528
$ws.Range("F16").Value = 1672
$ws.Range("F17").Value = 281
$ws.Range("F18").Value = 1213
$ws.Range("F19").Value = 211
$ws.Range("F20").Value = 472
$ws.Range("F22").Value = 114
$ws.Range("F23").Value = 125

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 4535
$ws.Range("F5").Value = 212
$ws.Range("F6").Value = 212
$ws.Range("F10").Value = 462

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5217
$ws.Range("F3").Value = 286

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 5217
$ws.Range("F4").Value = 286
$ws.Range("F5").Value = 27092
$ws.Range("G5").Value = "不可售"
$ws.Range("F6").Value = 4535
$ws.Range("F7").Value = 672
$ws.Range("F11").Value = 212
$ws.Range("F12").Value = 212
$ws.Range("F16").Value = 462
$ws.Range("F17").Value = 585
$ws.Range("F21").Value = 387
$ws.Range("F22").Value = 505
$ws.Range("F24").Value = 56
$ws.Range("F26").Value = 321
$ws.Range("F27").Value = 108
$ws.Range("F30").Value = 528
$ws.Range("F33").Value = 1672
$ws.Range("F34").Value = 281
$ws.Range("F35").Value = 1213
$ws.Range("F37").Value = 211
$ws.Range("F38").Value = 472
$ws.Range("F40").Value = 114
$ws.Range("F42").Value = 125
